# Generate Report for Handoff
# Adds a new handoff-status row (790af623-9f78-48ae-afa1-8243b2670b39) to the
# Overview / zh-cn / de-de sheets of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$newGuid = "790af623-9f78-48ae-afa1-8243b2670b39"
$newHash = "929dc35992db41f880409269abbfcf12085f3d6d"
$mdName  = "$newGuid.md"
$zhXlf   = "$newGuid.$newHash.zh-cn.xlf"
$deXlf   = "$newGuid.$newHash.de-de.xlf"

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/756731c8dfb711571ef2949eef68b03f1c112070/e2e/$mdName"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f96f76958731798a6dc18875c94281ea2e29c2d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bbeeb61e811aff4d6c470a6bd4a3b0bf16f28012/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" -> new row 3
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Cells.Item(3, 1).Value = $mdName
$ws1.Cells.Item(3, 2).Value = "Ready for handoff"
$ws1.Cells.Item(3, 3).Value = "Ready for handoff"
$ws1.Cells.Item(3, 4).Value = "2016-29-19 00:29:14"

$ws1.Hyperlinks.Add($ws1.Cells.Item(3, 1), $mdUrl, $null, $null, $mdName)

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn" -> new row 3
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Cells.Item(3, 1).Value = $mdName
$ws2.Cells.Item(3, 2).Value = ".md"
$ws2.Cells.Item(3, 3).Value = "Ready for handoff"
$ws2.Cells.Item(3, 4).Value = $zhXlf
$ws2.Cells.Item(3, 5).Value = "2016-03-19 00:29:11"
$ws2.Cells.Item(3, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(3, 8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(3, 9).Value = "Include"

$ws2.Hyperlinks.Add($ws2.Cells.Item(3, 1), $mdUrl, $null, $null, $mdName)
$ws2.Hyperlinks.Add($ws2.Cells.Item(3, 2), $mdUrl, $null, $null, ".md")
$ws2.Hyperlinks.Add($ws2.Cells.Item(3, 4), $zhUrl, $null, $null, $zhXlf)

# ---------------------------------------------------------------------------
# Sheet 3: "de-de" -> new row 3
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Cells.Item(3, 1).Value = $mdName
$ws3.Cells.Item(3, 2).Value = ".md"
$ws3.Cells.Item(3, 3).Value = "Ready for handoff"
$ws3.Cells.Item(3, 4).Value = $deXlf
$ws3.Cells.Item(3, 5).Value = "2016-03-19 00:29:14"
$ws3.Cells.Item(3, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(3, 8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(3, 9).Value = "Include"

$ws3.Hyperlinks.Add($ws3.Cells.Item(3, 1), $mdUrl, $null, $null, $mdName)
$ws3.Hyperlinks.Add($ws3.Cells.Item(3, 2), $mdUrl, $null, $null, ".md")
$ws3.Hyperlinks.Add($ws3.Cells.Item(3, 4), $deUrl, $null, $null, $deXlf)
